# Update FFXIV Leves profit figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 38596.258
$ws.Range("I135").Value = 43198.457
$ws.Range("J135").Value = 1778.6666
$ws.Range("K135").Value = 388786.113
$ws.Range("L135").Value = 16007.9994
$ws.Range("M135").Value = -386251.113
$ws.Range("N135").Value = -21077.9994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1188.2413
$ws.Range("I137").Value = 1011
$ws.Range("J137").Value = 1525
$ws.Range("K137").Value = 3033
$ws.Range("L137").Value = 4575
$ws.Range("M137").Value = -483
$ws.Range("N137").Value = -9675

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3354.742
$ws.Range("I138").Value = 2459.3333
$ws.Range("J138").Value = 4194.1875
$ws.Range("K138").Value = 7377.999899999999
$ws.Range("L138").Value = 12582.5625
$ws.Range("M138").Value = -2237.999899999999
$ws.Range("N138").Value = -22862.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2658.5
$ws.Range("I31").Value = 2658.5
$ws.Range("K31").Value = 2658.5
$ws.Range("M31").Value = -2364.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24255.656
$ws.Range("I32").Value = 4417.654
$ws.Range("J32").Value = 110220.336
$ws.Range("K32").Value = 4417.654
$ws.Range("L32").Value = 110220.336
$ws.Range("M32").Value = -4130.654
$ws.Range("N32").Value = -110794.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3007.3809
$ws.Range("I61").Value = 2842.9285
$ws.Range("J61").Value = 3336.2856
$ws.Range("K61").Value = 2842.9285
$ws.Range("L61").Value = 3336.2856
$ws.Range("M61").Value = -2630.9285
$ws.Range("N61").Value = -3760.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1279.8096
$ws.Range("I74").Value = 1165.3334
$ws.Range("J74").Value = 1966.6666
$ws.Range("K74").Value = 1165.3334
$ws.Range("L74").Value = 1966.6666
$ws.Range("M74").Value = -291.3334
$ws.Range("N74").Value = -3714.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1279.8096
$ws.Range("I77").Value = 1165.3334
$ws.Range("J77").Value = 1966.6666
$ws.Range("K77").Value = 5826.666999999999
$ws.Range("L77").Value = 9833.333000000001
$ws.Range("M77").Value = -1458.666999999999
$ws.Range("N77").Value = -18569.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1733.3903
$ws.Range("I132").Value = 1228.5
$ws.Range("J132").Value = 2318
$ws.Range("K132").Value = 3685.5
$ws.Range("L132").Value = 6954
$ws.Range("M132").Value = -1155.5
$ws.Range("N132").Value = -12014

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3007.3809
$ws.Range("I136").Value = 2842.9285
$ws.Range("J136").Value = 3336.2856
$ws.Range("K136").Value = 8528.7855
$ws.Range("L136").Value = 10008.8568
$ws.Range("M136").Value = -5978.7855
$ws.Range("N136").Value = -15108.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1297.0264
$ws.Range("I134").Value = 1064.1034
$ws.Range("J134").Value = 2047.5555
$ws.Range("K134").Value = 3192.3102
$ws.Range("L134").Value = 6142.666499999999
$ws.Range("M134").Value = -657.3101999999999
$ws.Range("N134").Value = -11212.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25643270
$ws.Range("I31").Value = 62501652
$ws.Range("J31").Value = 2653.5652
$ws.Range("K31").Value = 62501652
$ws.Range("L31").Value = 2653.5652
$ws.Range("M31").Value = -62501357
$ws.Range("N31").Value = -3243.5652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 25643270
$ws.Range("I34").Value = 62501652
$ws.Range("J34").Value = 2653.5652
$ws.Range("K34").Value = 62501652
$ws.Range("L34").Value = 2653.5652
$ws.Range("M34").Value = -62501450
$ws.Range("N34").Value = -3057.5652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 774.6129
$ws.Range("I107").Value = 747.43475
$ws.Range("J107").Value = 852.75
$ws.Range("K107").Value = 747.43475
$ws.Range("L107").Value = 852.75
$ws.Range("M107").Value = 1172.56525
$ws.Range("N107").Value = -4692.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2613.9473
$ws.Range("I132").Value = 1730.5834
$ws.Range("J132").Value = 4128.2856
$ws.Range("K132").Value = 5191.7502
$ws.Range("L132").Value = 12384.8568
$ws.Range("M132").Value = -2661.7502
$ws.Range("N132").Value = -17444.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9948.191999999999
$ws.Range("I134").Value = 7732.8945
$ws.Range("J134").Value = 15961.143
$ws.Range("K134").Value = 23198.6835
$ws.Range("L134").Value = 47883.429
$ws.Range("M134").Value = -20663.6835
$ws.Range("N134").Value = -52953.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 343.32144
$ws.Range("I5").Value = 271.91666
$ws.Range("J5").Value = 771.75
$ws.Range("K5").Value = 815.7499799999999
$ws.Range("L5").Value = 2315.25
$ws.Range("M5").Value = -703.7499799999999
$ws.Range("N5").Value = -2539.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 594.9530999999999
$ws.Range("I107").Value = 291.4634
$ws.Range("J107").Value = 1135.9565
$ws.Range("K107").Value = 874.3901999999999
$ws.Range("L107").Value = 3407.8695
$ws.Range("M107").Value = 1045.6098
$ws.Range("N107").Value = -7247.8695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 516.1667
$ws.Range("I113").Value = 350
$ws.Range("J113").Value = 549.4
$ws.Range("K113").Value = 1050
$ws.Range("L113").Value = 1648.2
$ws.Range("M113").Value = 1120
$ws.Range("N113").Value = -5988.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4385.5713
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 4783.1665
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 14349.4995
$ws.Range("M116").Value = -2558
$ws.Range("N116").Value = -21233.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3219
$ws.Range("I117").Value = 506.66666
$ws.Range("J117").Value = 4236.125
$ws.Range("K117").Value = 1519.99998
$ws.Range("L117").Value = 12708.375
$ws.Range("M117").Value = 1922.00002
$ws.Range("N117").Value = -19592.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 873.75555
$ws.Range("I131").Value = 325
$ws.Range("J131").Value = 992.4054
$ws.Range("K131").Value = 975
$ws.Range("L131").Value = 2977.2162
$ws.Range("M131").Value = 4065
$ws.Range("N131").Value = -13057.2162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 343.32144
$ws.Range("I135").Value = 271.91666
$ws.Range("J135").Value = 771.75
$ws.Range("K135").Value = 2447.24994
$ws.Range("L135").Value = 6945.75
$ws.Range("M135").Value = 87.7500600000003
$ws.Range("N135").Value = -12015.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3042.1035
$ws.Range("I132").Value = 2638.5715
$ws.Range("J132").Value = 4101.375
$ws.Range("K132").Value = 7915.7145
$ws.Range("L132").Value = 12304.125
$ws.Range("M132").Value = -5385.7145
$ws.Range("N132").Value = -17364.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3203.7273
$ws.Range("I122").Value = 2618.5454
$ws.Range("J122").Value = 3788.9092
$ws.Range("K122").Value = 7855.6362
$ws.Range("L122").Value = 11366.7276
$ws.Range("M122").Value = -5405.6362
$ws.Range("N122").Value = -16266.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2170.8
$ws.Range("I132").Value = 1968
$ws.Range("J132").Value = 2475
$ws.Range("K132").Value = 5904
$ws.Range("L132").Value = 7425
$ws.Range("M132").Value = -3374
$ws.Range("N132").Value = -12485

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1807.4445
$ws.Range("I136").Value = 1678.7693
$ws.Range("J136").Value = 2142
$ws.Range("K136").Value = 5036.3079
$ws.Range("L136").Value = 6426
$ws.Range("M136").Value = -2486.3079
$ws.Range("N136").Value = -11526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9640.200000000001
$ws.Range("I62").Value = 8280.4
$ws.Range("J62").Value = 11000
$ws.Range("K62").Value = 8280.4
$ws.Range("L62").Value = 11000
$ws.Range("M62").Value = -7656.4
$ws.Range("N62").Value = -12248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9640.200000000001
$ws.Range("I65").Value = 8280.4
$ws.Range("J65").Value = 11000
$ws.Range("K65").Value = 41402
$ws.Range("L65").Value = 55000
$ws.Range("M65").Value = -38282
$ws.Range("N65").Value = -61240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3664.7742
$ws.Range("I132").Value = 4059.1667
$ws.Range("J132").Value = 2312.5715
$ws.Range("K132").Value = 12177.5001
$ws.Range("L132").Value = 6937.7145
$ws.Range("M132").Value = -9647.500100000001
$ws.Range("N132").Value = -11997.7145
